# ML model retrained with all data
# Columns J and K (rows 1-51) hold the model's per-row threshold/weight
# outputs. After retraining on the full dataset every row now carries the
# same pair of constants: J = 0.6, K = 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1:J51").Value = 0.6
$ws.Range("K1:K51").Value = 1

# Selection moved to K1 (top of the newly-recomputed K column), matching
# where the author left the cursor after the refresh.
$ws.Range("K1:K51").Select()
